# Update the USD Amount (T2) and move the active selection from T4 to T3,
# matching the source workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# T2 changed from 543599 to 545611
$ws.Range("T2").Value = 545611

# Active selection moved from T4 to T3
$ws.Range("T3").Select()
